$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original sheet and add the three new sheets (S2, S3, S4)
#    in order, so they land after it in tab order / file order.
# ---------------------------------------------------------------------------
$tbl = $wb.Worksheets.Item(1)
$tbl.Name = "Table"

$s2 = $wb.Worksheets.Add($null, $tbl)
$s2.Name = "S2"

$s3 = $wb.Worksheets.Add($null, $s2)
$s3.Name = "S3"

$s4 = $wb.Worksheets.Add($null, $s3)
$s4.Name = "S4"

# ---------------------------------------------------------------------------
# 2. Fill in the raw survey scores on the "Table" sheet (rows 2-5).
# ---------------------------------------------------------------------------
$tbl.Range("A2").Value = 1
$tbl.Range("B2").Value = 1
$tbl.Range("C2").Value = 2
$tbl.Range("D2").Value = 2
$tbl.Range("F2").Value = 1
$tbl.Range("G2").Value = 1
$tbl.Range("H2").Value = 2
$tbl.Range("I2").Value = 1

$tbl.Range("A3").Value = 2
$tbl.Range("B3").Value = 2
$tbl.Range("D3").Value = 3
$tbl.Range("E3").Value = 0
$tbl.Range("F3").Value = 1
$tbl.Range("G3").Value = 0
$tbl.Range("H3").Value = 3
$tbl.Range("I3").Value = 1

$tbl.Range("A4").Value = 3
$tbl.Range("B4").Value = 2
$tbl.Range("C4").Value = 0
$tbl.Range("D4").Value = 0
$tbl.Range("E4").Value = 2
$tbl.Range("F4").Value = 3
$tbl.Range("G4").Value = 2
$tbl.Range("H4").Value = 0
$tbl.Range("I4").Value = 1

$tbl.Range("A5").Value = 4
$tbl.Range("G5").Value = 2
$tbl.Range("H5").Value = 1
$tbl.Range("I5").Value = 1

# ---------------------------------------------------------------------------
# 3. Populate the four question sheets with their section text.
#    Header cell gets a bold black Arial-BoldMT font, body cells get a
#    plain black ArialMT font. Set the font once on a template cell, then
#    propagate it with Copy / PasteSpecial(Formats) so we don't re-create
#    the font object for every single cell.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

# --- S2 : "Section 2" questionnaire -----------------------------------
$s2.Range("B4").Value = "Section 2. If I have a possible shortcoming in teamwork, it could be that:"
$s2.Range("B5").Value = "(a)I am not at ease unless meetings are well structured and controlled and generally well conducted."
$s2.Range("B6").Value = "(b)I am inclined to be too generous towards others who have a valid viewpoint that has not been given a proper airing."
$s2.Range("B7").Value = "(c)I have a tendency to talk a lot once the group gets on to new ideas."
$s2.Range("B8").Value = "(d)My objective outlook makes it difficult for me to join in readily and enthusiastically with colleagues."
$s2.Range("B9").Value = "(e)I am sometimes seen as forceful and authoritarian if there is a need to get something done."
$s2.Range("B10").Value = "(f) I find it difficult to lead from the front, perhaps because I am over-responsive to group atmosphere."
$s2.Range("B11").Value = "(g)I am apt to get too caught up in ideas that occur to me and so lose track of what is happening."
$s2.Range("B12").Value = "(h)My colleagues tend to see me as worrying unnecessarily over detail and the possibility that things may go wrong."

$s2.Range("B4").Font.Name = "Arial-BoldMT"
$s2.Range("B4").Font.Bold = $true
$s2.Range("B4").Font.Color = 0

$s2.Range("B5").Font.Name = "ArialMT"
$s2.Range("B5").Font.Color = 0

$s2.Range("B5").Copy()
$s2.Range("B6:B12").PasteSpecial($xlPasteFormats)
$s2.Range("B14").Font.Name = "ArialMT"
$s2.Range("B14").Font.Color = 0

# --- S3 : "Section 3" questionnaire -----------------------------------
$s3.Range("A1").Value = "Section 3. When involved in a project with other people:"
$s3.Range("A2").Value = "a"
$s3.Range("B2").Value = "(a)I have an aptitude for influencing people without pressurizing them."
$s3.Range("A3").Value = "b"
$s3.Range("B3").Value = "(b)My general vigilance percents careless mistakes and omissions being made."
$s3.Range("A4").Value = "c"
$s3.Range("B4").Value = "(c)I am ready to press for action to make sure that the meeting does not waste  time or lose sight of the main objective."
$s3.Range("A5").Value = "d"
$s3.Range("B5").Value = "(d)I can be counted on to contribute something original."
$s3.Range("A6").Value = "e"
$s3.Range("B6").Value = "(e)I am always ready to back a good suggestion in the common interest."
$s3.Range("A7").Value = "f"
$s3.Range("B7").Value = "(f) I am keen to look for the latest in new ideas and developments."
$s3.Range("A8").Value = "g"
$s3.Range("B8").Value = "(g)I believe my capacity for cool judgments is appreciated by others."
$s3.Range("A9").Value = "h"
$s3.Range("B9").Value = "(h)I can be relied upon to see that all essential work is organized."

$s3.Range("A1").Font.Name = "Arial-BoldMT"
$s3.Range("A1").Font.Bold = $true
$s3.Range("A1").Font.Color = 0

$s3.Range("B2").Font.Name = "ArialMT"
$s3.Range("B2").Font.Color = 0

$s3.Range("B2").Copy()
$s3.Range("B3:B9").PasteSpecial($xlPasteFormats)

# --- S4 : "Section 4" questionnaire -----------------------------------
$s4.Range("A1").Value = "Section 4. My characteristic approach to group work is that:"
$s4.Range("A2").Value = "a"
$s4.Range("B2").Value = "(a)I have a quiet interest in getting to know colleagues better."
$s4.Range("A3").Value = "b"
$s4.Range("B3").Value = "(b)I am not reluctant to challenge the views of others or to hold a minority view myself."
$s4.Range("A4").Value = "c"
$s4.Range("B4").Value = "(c)I can usually find a line of argument to refute unsound propositions."
$s4.Range("A5").Value = "d"
$s4.Range("B5").Value = "(d)I think I have a talent for making things work once a plan has to be put into operation."
$s4.Range("A6").Value = "e"
$s4.Range("B6").Value = "(e)I have a tendency to avoid the obvious and to come out with the unexpected."
$s4.Range("A7").Value = "f"
$s4.Range("B7").Value = "(f) I bring a touch of perfectionism to any team job I undertake."
$s4.Range("A8").Value = "g"
$s4.Range("B8").Value = "(g)I am ready to make use of contacts outside the group itself."
$s4.Range("A9").Value = "h"
$s4.Range("B9").Value = "(h)While I am interested in all views, I have no hesitation in making up my mind once a decision has to be made."

$s4.Range("A1").Font.Name = "Arial-BoldMT"
$s4.Range("A1").Font.Bold = $true
$s4.Range("A1").Font.Color = 0

$s4.Range("B2").Font.Name = "ArialMT"
$s4.Range("B2").Font.Color = 0

$s4.Range("B2").Copy()
$s4.Range("B3:B9").PasteSpecial($xlPasteFormats)
$s4.Range("B10").Font.Name = "ArialMT"
$s4.Range("B10").Font.Color = 0

$s4.PageSetup.PaperSize = 9
$s4.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4. Selections & active tab -- S2 is the tab that's on top when the file
#    is opened, each sheet keeps its own last-used selection.
# ---------------------------------------------------------------------------
$tbl.Range("E5").Select()
$s3.Range("C14").Select()
$s4.Range("G22").Select()
$s2.Range("C22").Select()
$s2.Activate()

$wb.Save()
